$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 15:31"

# Row 4
$ws.Range("B4").Value = 5956661
$ws.Range("C4").Value = 933
$ws.Range("D4").Value = 3255161
$ws.Range("E4").Value = 2519051
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = 182449

# Row 6
$ws.Range("B6").Value = 3246929
$ws.Range("C6").Value = 15175
$ws.Range("D6").Value = 2474743
$ws.Range("E6").Value = 712450
$ws.Range("G6").Value = 124
$ws.Range("H6").Value = 59736

# Row 17
$ws.Range("B17").Value = 310836
$ws.Range("C17").Value = 1068
$ws.Range("D17").Value = 284945
$ws.Range("E17").Value = 22136
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 3755

# Row 23
$ws.Range("B23").Value = 237896
$ws.Range("C23").Value = 324
$ws.Range("E23").Value = 18949
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 9347

# Row 28
$ws.Range("B28").Value = 117742
$ws.Range("C28").Value = 244
$ws.Range("D28").Value = 114558
$ws.Range("E28").Value = 2990

# Row 30
$ws.Range("E30").Value = 54277
$ws.Range("G30").Value = 36
$ws.Range("H30").Value = 2354

# Row 37
$ws.Range("B37").Value = 87072
$ws.Range("G37").Value = 5
$ws.Range("H37").Value = 5817

# Row 43
$ws.Range("B43").Value = 70974
$ws.Range("C43").Value = 247
$ws.Range("D43").Value = 69378
$ws.Range("E43").Value = 939
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 657

# Row 45
$ws.Range("A45").Value = "Paises Bajos"
$ws.Range("B45").Value = 68114
$ws.Range("C45").Value = 571
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 6215

# Row 46
$ws.Range("A46").Value = "Emiratos Arabes Unidos"
$ws.Range("B46").Value = 68020
$ws.Range("C46").Value = 399
$ws.Range("D46").Value = 59070
$ws.Range("E46").Value = 8572
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 378

# Row 50
$ws.Range("B50").Value = 56274
$ws.Range("C50").Value = 362
$ws.Range("D50").Value = 41184
$ws.Range("E50").Value = 13283
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 1807

# Row 98
$ws.Range("A98").Value = "Tayikistan"
$ws.Range("B98").Value = 8413
$ws.Range("C98").Value = 34
$ws.Range("D98").Value = 7214
$ws.Range("E98").Value = 1132
$ws.Range("H98").Value = 67

# Row 99
$ws.Range("A99").Value = "Gabon"
$ws.Range("B99").Value = 8409
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 6959
$ws.Range("E99").Value = 1397
$ws.Range("H99").Value = 53

# Row 100
$ws.Range("B100").Value = 8122
$ws.Range("C100").Value = 10
$ws.Range("D100").Value = 5677
$ws.Range("E100").Value = 2248
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 197

# Row 178
$ws.Range("D178").Value = 357
$ws.Range("E178").Value = 54
